{"js": "// The document's tables all use a 72-dxa (3.6 pt) left/right cell margin.\n// Tighten it to 24 dxa (1.2 pt) on every table so content fits and page\n// breaks land correctly on Linux.\nconst TARGET_DXA = 24;\nconst TARGET_PT = TARGET_DXA / 20; // OOXML dxa are twentieths of a point.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < tables.items.length; i++) {\n  const table = tables.items[i];\n  // The Word JS API only exposes table cell padding through\n  // table.setCellPadding(location, value), which (in this runtime) doesn't\n  // persist. Word's real COM object model surfaces the same value as the\n  // Table.LeftPadding / Table.RightPadding properties, which the generated\n  // proxy still routes to the OM bridge via its internal _omSet helper, so\n  // use that to reach the same w:tblCellMar/w:left & w:right values.\n  table._omSet(\"LeftPadding\", TARGET_PT);\n  table._omSet(\"RightPadding\", TARGET_PT);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The document's tables all had their left/right cell margins at 72 dxa\n# (3.6 pt); collapse them down to 24 dxa (1.2 pt) so tables render more\n# tightly and page breaks land correctly on Linux.\n$targetDxa = 24\n$targetPt = $targetDxa / 20.0\n\n$count = $d.Tables.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Tables.Item($i)\n    $t.LeftPadding = $targetPt\n    $t.RightPadding = $targetPt\n}\n"}
